$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.150.23"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.822.39"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6016"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07054"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.97%  "
$ws.Range("E9").Value = "  -3.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07641"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").Value = "1.815.88"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.779"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6279"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009901"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.93%  "
$ws.Range("D16").Value = "2.064.87"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("E17").Value = "  -4.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.833"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "29.141.64"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "225.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.963"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.992"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1299"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("E28").Value = "  -4.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.486"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06224"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -15.60%  "
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.822"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.787"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.39%  "
$ws.Range("E34").Value = "  -1.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.734"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6361"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.541"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("D38").Value = "1.211.99"
$ws.Range("E38").Value = "  -1.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.720"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01733"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.474"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9008"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "1.978.67"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000117"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.475"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.22%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.587"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4552"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05501"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.85%  "
